# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2404"
#   "<name>_new" -> "<name>_FV2410"
# Then turn the used range into an Excel Table (ListObject) and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) -------------------------------------

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2404"
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $fv2404Headers[$i]
}

# Column K (11) is "diff" and stays untouched.

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2410"
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $fv2410Headers[$i]
}

# --- 2. Turn the used range into a Table (ListObject) ----------------------

$usedRange = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row -----------------------------------------------

$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
